$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M18").Value = -206.3
$ws.Range("K18").Value = 490.3
$ws.Range("I18").Value = 490.3
$ws.Range("H18").Value = 490.3
$ws.Range("I96").Value = 2195.9524
$ws.Range("H96").Value = 2739.3333
$ws.Range("L96").Value = 13923.4995
$ws.Range("K96").Value = 6587.8572
$ws.Range("N96").Value = -16669.4995
$ws.Range("M96").Value = -5214.8572
$ws.Range("J96").Value = 4641.1665
$ws.Range("L103").Value = 2549.25
$ws.Range("H103").Value = 824.5714
$ws.Range("N103").Value = -3721.25
$ws.Range("I103").Value = 791
$ws.Range("J103").Value = 849.75
$ws.Range("M103").Value = -1787
$ws.Range("K103").Value = 2373
$ws.Range("I135").Value = 1912.5
$ws.Range("K135").Value = 17212.5
$ws.Range("M135").Value = -14677.5
$ws.Range("H135").Value = 22729348
$ws.Range("J136").Value = 79444.22
$ws.Range("H136").Value = 79444.22
$ws.Range("L136").Value = 79444.22
$ws.Range("N136").Value = -89644.22
$ws.Range("I137").Value = 6860
$ws.Range("K137").Value = 20580
$ws.Range("J137").Value = 2213.5715
$ws.Range("L137").Value = 6640.7145
$ws.Range("M137").Value = -18030
$ws.Range("H137").Value = 4149.5835
$ws.Range("N137").Value = -11740.7145
$ws.Range("J138").Value = 10422788
$ws.Range("L138").Value = 31268364
$ws.Range("I138").Value = 2139.1538
$ws.Range("K138").Value = 6417.4614
$ws.Range("H138").Value = 7412378.5
$ws.Range("M138").Value = -1277.4614
$ws.Range("N138").Value = -31278644
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L44").Value = 39495
$ws.Range("H44").Value = 39495
$ws.Range("J44").Value = 39495
$ws.Range("N44").Value = -40471
$ws.Range("N54").Value = -33785.5
$ws.Range("J54").Value = 32247.5
$ws.Range("H54").Value = 32247.5
$ws.Range("L54").Value = 32247.5
$ws.Range("H61").Value = 43482990
$ws.Range("I61").Value = 62502864
$ws.Range("M61").Value = -62502652
$ws.Range("K61").Value = 62502864
$ws.Range("I74").Value = 55619132
$ws.Range("M74").Value = -55618258
$ws.Range("J74").Value = 3573.818
$ws.Range("H74").Value = 34523576
$ws.Range("L74").Value = 3573.818
$ws.Range("N74").Value = -5321.818
$ws.Range("K74").Value = 55619132
$ws.Range("I77").Value = 55619132
$ws.Range("H77").Value = 34523576
$ws.Range("L77").Value = 17869.09
$ws.Range("N77").Value = -26605.09
$ws.Range("M77").Value = -278091292
$ws.Range("K77").Value = 278095660
$ws.Range("J77").Value = 3573.818
$ws.Range("I136").Value = 62502864
$ws.Range("K136").Value = 187508592
$ws.Range("H136").Value = 43482990
$ws.Range("M136").Value = -187506042
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L11").Value = 785
$ws.Range("H11").Value = 659.2222
$ws.Range("K11").Value = 502
$ws.Range("J11").Value = 785
$ws.Range("I11").Value = 502
$ws.Range("M11").Value = -362
$ws.Range("N11").Value = -1065
$ws.Range("I86").Value = 10671.615
$ws.Range("H86").Value = 16652.521
$ws.Range("J86").Value = 24427.7
$ws.Range("K86").Value = 10671.615
$ws.Range("L86").Value = 24427.7
$ws.Range("M86").Value = -9548.615
$ws.Range("N86").Value = -26673.7
$ws.Range("N89").Value = -133370.5
$ws.Range("I89").Value = 10671.615
$ws.Range("H89").Value = 16652.521
$ws.Range("L89").Value = 122138.5
$ws.Range("J89").Value = 24427.7
$ws.Range("K89").Value = 53358.075
$ws.Range("M89").Value = -47742.075
$ws.Range("N99").Value = -8598
$ws.Range("K99").Value = 1688.0834
$ws.Range("H99").Value = 2839.2354
$ws.Range("J99").Value = 5602
$ws.Range("I99").Value = 1688.0834
$ws.Range("M99").Value = -190.0834
$ws.Range("L99").Value = 5602
$ws.Range("H107").Value = 3287.3157
$ws.Range("N107").Value = -9328
$ws.Range("L107").Value = 5488
$ws.Range("I107").Value = 2142.96
$ws.Range("M107").Value = -222.96
$ws.Range("K107").Value = 2142.96
$ws.Range("J107").Value = 5488
$ws.Range("M134").Value = -7326
$ws.Range("I134").Value = 3287
$ws.Range("H134").Value = 3249.9355
$ws.Range("K134").Value = 9861
$ws.Range("H141").Value = 20001
$ws.Range("K141").Value = 20001
$ws.Range("I141").Value = 20001
$ws.Range("M141").Value = -14821
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N31").Value = -50007310
$ws.Range("H31").Value = 26321628
$ws.Range("L31").Value = 50006720
$ws.Range("J31").Value = 50006720
$ws.Range("N34").Value = -50007124
$ws.Range("H34").Value = 26321628
$ws.Range("J34").Value = 50006720
$ws.Range("L34").Value = 50006720
$ws.Range("H52").Value = 148967.75
$ws.Range("J52").Value = 148967.75
$ws.Range("N52").Value = -149555.75
$ws.Range("L52").Value = 148967.75
$ws.Range("N64").Value = -40496
$ws.Range("L64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("H64").Value = 40000
$ws.Range("N67").Value = -41716
$ws.Range("J67").Value = 40000
$ws.Range("H67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("J94").Value = 2194
$ws.Range("H94").Value = 2068.2727
$ws.Range("L94").Value = 2194
$ws.Range("N94").Value = -3096
$ws.Range("N105").Value = -22682.666
$ws.Range("K105").Value = 1516.5
$ws.Range("L105").Value = 19188.666
$ws.Range("H105").Value = 10352.583
$ws.Range("M105").Value = 230.5
$ws.Range("I105").Value = 1516.5
$ws.Range("J105").Value = 19188.666
$ws.Range("L131").Value = 76399.336
$ws.Range("H131").Value = 76399.336
$ws.Range("N131").Value = -86479.336
$ws.Range("J131").Value = 76399.336
$ws.Range("M134").Value = -2179.125
$ws.Range("I134").Value = 1571.375
$ws.Range("H134").Value = 1741.9
$ws.Range("L134").Value = 7272
$ws.Range("K134").Value = 4714.125
$ws.Range("J134").Value = 2424
$ws.Range("N134").Value = -12342
$ws.Range("H141").Value = 241871
$ws.Range("L141").Value = 342141.28
$ws.Range("N141").Value = -352501.28
$ws.Range("J141").Value = 342141.28
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N37").Value = -600216.5
$ws.Range("L37").Value = 599992.5
$ws.Range("J37").Value = 199997.5
$ws.Range("H37").Value = 199997.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3318
$ws.Range("L97").Value = 4138.7144
$ws.Range("M97").Value = -907
$ws.Range("N97").Value = -5130.7144
$ws.Range("K97").Value = 1403
$ws.Range("J97").Value = 4138.7144
$ws.Range("I97").Value = 1403
$ws.Range("I102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H102").Value = 4835.6665
$ws.Range("K102").Value = 0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 7017.5
$ws.Range("K31").Value = 228
$ws.Range("M31").Value = 20
$ws.Range("I31").Value = 228
$ws.Range("I40").Value = 3191.353
$ws.Range("L40").Value = 5876.1816
$ws.Range("H40").Value = 3847.6445
$ws.Range("J40").Value = 5876.1816
$ws.Range("M40").Value = -3055.353
$ws.Range("N40").Value = -6148.1816
$ws.Range("K40").Value = 3191.353
$ws.Range("I46").Value = 788.08
$ws.Range("M46").Value = -600.08
$ws.Range("K46").Value = 788.08
$ws.Range("H46").Value = 1354.9032
$ws.Range("H55").Value = 942.2308
$ws.Range("K55").Value = 468.625
$ws.Range("I55").Value = 468.625
$ws.Range("M55").Value = -295.625
$ws.Range("I100").Value = 2071.8572
$ws.Range("K100").Value = 2071.8572
$ws.Range("H100").Value = 2696.4167
$ws.Range("M100").Value = -1530.8572
$ws.Range("I122").Value = 3696.8948
$ws.Range("N122").Value = -21828.1432
$ws.Range("K122").Value = 11090.6844
$ws.Range("H122").Value = 4220.769
$ws.Range("J122").Value = 5642.7144
$ws.Range("L122").Value = 16928.1432
$ws.Range("M122").Value = -8640.6844
$ws.Range("M132").Value = -2169.9998
$ws.Range("K132").Value = 4699.9998
$ws.Range("H132").Value = 250001710
$ws.Range("J132").Value = 400001800
$ws.Range("N132").Value = -1200010460
$ws.Range("L132").Value = 1200005400
$ws.Range("I132").Value = 1566.6666
$ws.Range("I136").Value = 1978.2391
$ws.Range("K136").Value = 5934.7173
$ws.Range("J136").Value = 1666.6666
$ws.Range("H136").Value = 1959.1632
$ws.Range("L136").Value = 4999.9998
$ws.Range("N136").Value = -10099.9998
$ws.Range("M136").Value = -3384.7173
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L131").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("J131").Value = 0
